$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.339.62"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.275.03"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "308.80"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "97.78"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").Value = "0.0819"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "6.84"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "2.627.61"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "2.272.06"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "42.225.59"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "12.32"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "67.71"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "236.95"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "23.61"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "37.46"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "9.57"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "163.67"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").Value = "0.0733"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "4.18"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E42").Value = "  -6.04%  "
$ws.Range("D43").Value = "1.946.56"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").Value = "0.0284"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "18.76"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").Value = "9.80"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "54.27"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "2.499.49"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "92.21"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "71.60"
$ws.Range("E51").Value = "  -1.62%  "
